$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "NeedTar" column (column I) and
# give it the header "AutoAtkDis", matching the new property added for skills.
$ws.Columns("I:I").Insert()
$ws.Range("I1").Value = "AutoAtkDis"

# Fill the new AutoAtkDis column with a value of 1 for every data row (2-9).
$ws.Range("I2:I9").Value = 1

# Reflect the selection the author left the sheet with after adding the column.
$ws.Range("I2:I9").Select() | Out-Null
